# Updates the cryptos price/volume snapshot (Sheet1, columns D = Price,
# E = Volume(1h)) to the new scraped values, per the commit diff.
#
# Values that are ambiguous with genuine numbers (e.g. "584.42", "0.596")
# are written with a leading apostrophe so Excel stores them as text,
# exactly matching the original inline-string/text representation of the
# Price column (which mixes plain numbers, dotted "thousands" strings like
# "62.958.73", and plain decimals like "584.42" - all of them as text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Formula  = "62.958.73"
$ws.Cells.Item(3, 4).Formula  = "2.590.02"
$ws.Cells.Item(3, 5).Formula  = "  +1.51%  "
$ws.Cells.Item(4, 5).Formula  = "  -0.04%  "
$ws.Cells.Item(5, 4).Formula  = "'584.42"
$ws.Cells.Item(5, 5).Formula  = "  +0.01%  "
$ws.Cells.Item(6, 4).Formula  = "'147.03"
$ws.Cells.Item(6, 5).Formula  = "  -0.17%  "
$ws.Cells.Item(7, 5).Formula  = "  +0.00%  "
$ws.Cells.Item(8, 4).Formula  = "'0.596"
$ws.Cells.Item(8, 5).Formula  = "  +2.19%  "
$ws.Cells.Item(9, 5).Formula  = "  +1.80%  "
$ws.Cells.Item(10, 5).Formula = "  +1.99%  "
$ws.Cells.Item(11, 5).Formula = "  -0.07%  "
$ws.Cells.Item(12, 5).Formula = "  -0.42%  "
$ws.Cells.Item(13, 4).Formula = "'27.29"
$ws.Cells.Item(13, 5).Formula = "  -0.28%  "
$ws.Cells.Item(14, 4).Formula = "3.053.90"
$ws.Cells.Item(14, 5).Formula = "  +1.53%  "
$ws.Cells.Item(15, 4).Formula = "62.872.63"
$ws.Cells.Item(15, 5).Formula = "  -0.13%  "
$ws.Cells.Item(16, 5).Formula = "  +2.71%  "
$ws.Cells.Item(17, 4).Formula = "2.586.10"
$ws.Cells.Item(17, 5).Formula = "  +1.55%  "
$ws.Cells.Item(18, 5).Formula = "  -0.55%  "
$ws.Cells.Item(19, 4).Formula = "'341.99"
$ws.Cells.Item(19, 5).Formula = "  +1.65%  "
$ws.Cells.Item(20, 5).Formula = "  +1.31%  "
$ws.Cells.Item(21, 5).Formula = "  -1.24%  "
$ws.Cells.Item(22, 5).Formula = "  -0.03%  "
$ws.Cells.Item(23, 4).Formula = "'67.25"
$ws.Cells.Item(23, 5).Formula = "  +2.20%  "
$ws.Cells.Item(24, 4).Formula = "2.713.88"
$ws.Cells.Item(24, 5).Formula = "  +1.73%  "
$ws.Cells.Item(25, 5).Formula = "  -1.56%  "
$ws.Cells.Item(26, 5).Formula = "  -1.48%  "
$ws.Cells.Item(27, 5).Formula = "  -0.11%  "
$ws.Cells.Item(28, 4).Formula = "'8.34"
$ws.Cells.Item(28, 5).Formula = "  -0.49%  "
$ws.Cells.Item(29, 4).Formula = "'7.85"
$ws.Cells.Item(29, 5).Formula = "  +5.73%  "
$ws.Cells.Item(30, 5).Formula = "  -1.32%  "
$ws.Cells.Item(31, 4).Formula = "'1.94"
$ws.Cells.Item(31, 5).Formula = "  +0.70%  "
$ws.Cells.Item(32, 4).Formula = "'476.68"
$ws.Cells.Item(32, 5).Formula = "  +14.44%  "
$ws.Cells.Item(33, 5).Formula = "  +1.09%  "
$ws.Cells.Item(34, 4).Formula = "'176.87"
$ws.Cells.Item(34, 5).Formula = "  -0.51%  "
$ws.Cells.Item(35, 5).Formula = "  +4.21%  "
$ws.Cells.Item(36, 5).Formula = "  +0.02%  "
$ws.Cells.Item(37, 4).Formula = "'0.405"
$ws.Cells.Item(37, 5).Formula = "  +1.19%  "
$ws.Cells.Item(38, 4).Formula = "'19.04"
$ws.Cells.Item(38, 5).Formula = "  -0.59%  "
$ws.Cells.Item(39, 5).Formula = "  +3.96%  "
$ws.Cells.Item(40, 5).Formula = "  +0.00%  "
$ws.Cells.Item(41, 5).Formula = "  -2.37%  "
$ws.Cells.Item(42, 4).Formula = "'158.35"
$ws.Cells.Item(42, 5).Formula = "  +4.77%  "
$ws.Cells.Item(43, 5).Formula = "  -0.07%  "
$ws.Cells.Item(44, 4).Formula = "'21.33"
$ws.Cells.Item(44, 5).Formula = "  +2.13%  "
$ws.Cells.Item(45, 5).Formula = "  +5.34%  "
$ws.Cells.Item(46, 5).Formula = "  +0.44%  "
$ws.Cells.Item(47, 5).Formula = "  -0.11%  "
$ws.Cells.Item(48, 4).Formula = "'0.0237"
$ws.Cells.Item(48, 5).Formula = "  -1.02%  "
$ws.Cells.Item(49, 5).Formula = "  +0.13%  "
$ws.Cells.Item(50, 5).Formula = "  +0.73%  "
$ws.Cells.Item(51, 5).Formula = "  +1.06%  "
